$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 3 (summary row): turn plain values into formulas
# ---------------------------------------------------------------------
$ws.Range("D3").Formula = '=SUM(D4:D34)'
$ws.Range("E3").Formula = '=SUM(F3:I3)'
$ws.Range("F3").Formula = '=SUM(F4:F31)'
$ws.Range("G3").Formula = '=SUM(G4:G31)'
$ws.Range("H3").Formula = '=SUM(H4:H31)'
$ws.Range("I3").Formula = '=SUM(I4:I31)'
$ws.Range("J3").Formula = '=SUM(J4:J31)'
$ws.Range("K3").Formula = '=((F3*5)+(G3*4)+(H3*3)+(I3*2))/E3'
$ws.Range("L3").Formula = '=((F3+G3)/D3)'
$ws.Range("M3").Formula = '=(F3+G3+H3)/D3'

# ---------------------------------------------------------------------
# 2. Student rows 4-29: add "кол-во" (D) and grade-indicator (F:J) columns
# ---------------------------------------------------------------------
for ($r = 4; $r -le 29; $r++) {
    $ws.Range("D$r").Value = 1
}

# Row 4 uses its own (non-shared) formula for F4, then a shared formula
# group for G4:J4.
$ws.Range("F4").Formula = '=IF(F$2=$C4,1,0)'
$ws.Range("G4:J4").Formula = '=IF(G$2=$C4,1,0)'

# Rows 5-29 share one formula group across F:J.
$ws.Range("F5:J29").Formula = '=IF(F$2=$C5,1,0)'

# ---------------------------------------------------------------------
# 3. Conditional formatting: 3-color scale over F4:J29
# ---------------------------------------------------------------------
$rng = $ws.Range("F4:J29")
$cf = $rng.FormatConditions.AddColorScale(3)
$cf.ColorScaleCriteria.Item(1).Type = 1
$cf.ColorScaleCriteria.Item(1).FormatColor.Color = 7039480
$cf.ColorScaleCriteria.Item(2).Type = 4
$cf.ColorScaleCriteria.Item(2).Value = 50
$cf.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167
$cf.ColorScaleCriteria.Item(3).Type = 2
$cf.ColorScaleCriteria.Item(3).FormatColor.Color = 8109667

# ---------------------------------------------------------------------
# 4. Recalculate and fix up the active selection on the frozen pane
# ---------------------------------------------------------------------
$excel.Calculate()
$ws.Range("D3:M3").Select() | Out-Null
